$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value  = "done. added note in ch4 TSP domain description.  Added discussion on continuous domains in future work"
$ws.Range("C17").Value = "done.  Added section to future work"
$ws.Range("D38").Value = "TSP and knapsack are made discrete for this work.  Elevator is a discrete domain.  Adversarial problems:  Non contiguous solution spaces could fool SBE-trace"
$ws.Range("C16").Value = "done. Added description in sec 5.2"
$ws.Range("C35").Value = "done. added note in ch4 TSP domain description.  Added discussion on continuous domains in future work"

$ws.Range("C8").WrapText = $true
$ws.Range("C16").WrapText = $true
$ws.Range("C17").WrapText = $true
$ws.Range("C35").WrapText = $true
$ws.Range("D38").WrapText = $true

$excel.ActiveWindow.Zoom = 75
$ws.Range("C8").Select()
